$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column CN (92) needs the same 12-char width as the rest of the table ---
$ws.Columns(92).ColumnWidth = 11.17

# --- Row 1: new date header, CN1 ("2024/12/09") as literal text, style matches CM1 (s=1) ---
$ws.Range("CN1").NumberFormat = "@"
$ws.Range("CN1").Value = "2024/12/09"
$ws.Range("A2").Copy()
$ws.Range("CN1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 2-53: new numeric data for 2024/12/09, with fill matching value thresholds ---
$ws.Range("CN2").Value = 138
$ws.Range("N2").Copy()
$ws.Range("CN2").PasteSpecial(-4122)
$ws.Range("CN3").Value = 266.4
$ws.Range("A2").Copy()
$ws.Range("CN3").PasteSpecial(-4122)
$ws.Range("CN4").Value = 188.6
$ws.Range("A2").Copy()
$ws.Range("CN4").PasteSpecial(-4122)
$ws.Range("CN5").Value = 170.8
$ws.Range("A2").Copy()
$ws.Range("CN5").PasteSpecial(-4122)
$ws.Range("CN6").Value = 197.6
$ws.Range("A2").Copy()
$ws.Range("CN6").PasteSpecial(-4122)
$ws.Range("CN7").Value = 143.5
$ws.Range("A2").Copy()
$ws.Range("CN7").PasteSpecial(-4122)
$ws.Range("CN8").Value = 266.7
$ws.Range("A2").Copy()
$ws.Range("CN8").PasteSpecial(-4122)
$ws.Range("CN9").Value = 232.2
$ws.Range("A2").Copy()
$ws.Range("CN9").PasteSpecial(-4122)
$ws.Range("CN10").Value = 136.9
$ws.Range("N2").Copy()
$ws.Range("CN10").PasteSpecial(-4122)
$ws.Range("CN11").Value = 165.3
$ws.Range("A2").Copy()
$ws.Range("CN11").PasteSpecial(-4122)
$ws.Range("CN12").Value = 186.1
$ws.Range("A2").Copy()
$ws.Range("CN12").PasteSpecial(-4122)
$ws.Range("CN13").Value = 159.2
$ws.Range("A2").Copy()
$ws.Range("CN13").PasteSpecial(-4122)
$ws.Range("CN14").Value = 122.6
$ws.Range("D2").Copy()
$ws.Range("CN14").PasteSpecial(-4122)
$ws.Range("CN15").Value = 146.4
$ws.Range("A2").Copy()
$ws.Range("CN15").PasteSpecial(-4122)
$ws.Range("CN16").Value = 192.5
$ws.Range("A2").Copy()
$ws.Range("CN16").PasteSpecial(-4122)
$ws.Range("CN17").Value = 163.8
$ws.Range("A2").Copy()
$ws.Range("CN17").PasteSpecial(-4122)
$ws.Range("CN18").Value = 153.9
$ws.Range("A2").Copy()
$ws.Range("CN18").PasteSpecial(-4122)
$ws.Range("CN19").Value = 136.3
$ws.Range("N2").Copy()
$ws.Range("CN19").PasteSpecial(-4122)
$ws.Range("CN20").Value = 178
$ws.Range("A2").Copy()
$ws.Range("CN20").PasteSpecial(-4122)
$ws.Range("CN21").Value = 161.3
$ws.Range("A2").Copy()
$ws.Range("CN21").PasteSpecial(-4122)
$ws.Range("CN22").Value = 257
$ws.Range("A2").Copy()
$ws.Range("CN22").PasteSpecial(-4122)
$ws.Range("CN23").Value = 134.6
$ws.Range("N2").Copy()
$ws.Range("CN23").PasteSpecial(-4122)
$ws.Range("CN24").Value = 176.4
$ws.Range("A2").Copy()
$ws.Range("CN24").PasteSpecial(-4122)
$ws.Range("CN25").Value = 124.6
$ws.Range("D2").Copy()
$ws.Range("CN25").PasteSpecial(-4122)
$ws.Range("CN26").Value = 164.6
$ws.Range("A2").Copy()
$ws.Range("CN26").PasteSpecial(-4122)
$ws.Range("CN27").Value = 148.1
$ws.Range("A2").Copy()
$ws.Range("CN27").PasteSpecial(-4122)
$ws.Range("CN28").Value = 337.7
$ws.Range("A2").Copy()
$ws.Range("CN28").PasteSpecial(-4122)
$ws.Range("CN29").Value = 184.3
$ws.Range("A2").Copy()
$ws.Range("CN29").PasteSpecial(-4122)
$ws.Range("CN30").Value = 220.1
$ws.Range("A2").Copy()
$ws.Range("CN30").PasteSpecial(-4122)
$ws.Range("CN31").Value = 125.1
$ws.Range("N2").Copy()
$ws.Range("CN31").PasteSpecial(-4122)
$ws.Range("CN32").Value = 142.1
$ws.Range("A2").Copy()
$ws.Range("CN32").PasteSpecial(-4122)
$ws.Range("CN33").Value = 217.3
$ws.Range("A2").Copy()
$ws.Range("CN33").PasteSpecial(-4122)
$ws.Range("CN34").Value = 130.7
$ws.Range("N2").Copy()
$ws.Range("CN34").PasteSpecial(-4122)
$ws.Range("CN35").Value = 150.5
$ws.Range("A2").Copy()
$ws.Range("CN35").PasteSpecial(-4122)
$ws.Range("CN36").Value = 133.4
$ws.Range("N2").Copy()
$ws.Range("CN36").PasteSpecial(-4122)
$ws.Range("CN37").Value = 124.4
$ws.Range("D2").Copy()
$ws.Range("CN37").PasteSpecial(-4122)
$ws.Range("CN38").Value = 221.2
$ws.Range("A2").Copy()
$ws.Range("CN38").PasteSpecial(-4122)
$ws.Range("CN39").Value = 159.4
$ws.Range("A2").Copy()
$ws.Range("CN39").PasteSpecial(-4122)
$ws.Range("CN40").Value = 119.8
$ws.Range("D2").Copy()
$ws.Range("CN40").PasteSpecial(-4122)
$ws.Range("CN41").Value = 145.6
$ws.Range("A2").Copy()
$ws.Range("CN41").PasteSpecial(-4122)
$ws.Range("CN42").Value = 147.2
$ws.Range("A2").Copy()
$ws.Range("CN42").PasteSpecial(-4122)
$ws.Range("CN43").Value = 163.8
$ws.Range("A2").Copy()
$ws.Range("CN43").PasteSpecial(-4122)
$ws.Range("CN44").Value = 150.4
$ws.Range("A2").Copy()
$ws.Range("CN44").PasteSpecial(-4122)
$ws.Range("CN45").Value = 136.6
$ws.Range("N2").Copy()
$ws.Range("CN45").PasteSpecial(-4122)
$ws.Range("CN46").Value = 171.8
$ws.Range("A2").Copy()
$ws.Range("CN46").PasteSpecial(-4122)
$ws.Range("CN47").Value = 151.7
$ws.Range("A2").Copy()
$ws.Range("CN47").PasteSpecial(-4122)
$ws.Range("CN48").Value = 218.6
$ws.Range("A2").Copy()
$ws.Range("CN48").PasteSpecial(-4122)
$ws.Range("CN49").Value = 163.3
$ws.Range("A2").Copy()
$ws.Range("CN49").PasteSpecial(-4122)
$ws.Range("CN50").Value = 159.9
$ws.Range("A2").Copy()
$ws.Range("CN50").PasteSpecial(-4122)
$ws.Range("CN51").Value = 161.2
$ws.Range("A2").Copy()
$ws.Range("CN51").PasteSpecial(-4122)
$ws.Range("CN52").Value = 145.4
$ws.Range("A2").Copy()
$ws.Range("CN52").PasteSpecial(-4122)
$ws.Range("CN53").Value = 143.9
$ws.Range("A2").Copy()
$ws.Range("CN53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
